$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("N13")
$rng.Value = 3000
try {
  $rng.BorderAround(1, -4138)
  Write-Host "BorderAround OK"
} catch {
  Write-Host "ERR: $_"
}
